# Updates cryptos list values (price + 1h volume change) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
#
# Cells are plain text in the sheet (t="inlineStr"/shared string), not numbers.
# For D-column values that are purely numeric-looking (e.g. "560.34"), a
# leading apostrophe is used so Excel stores them as TEXT (matching the
# original formatting) instead of silently converting them to a float -
# exactly what a leading apostrophe does when typed into Excel by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.902.12'
$ws.Range('E2').Value = '  +0.17%  '

$ws.Range('D3').Value = '2.439.35'
$ws.Range('E3').Value = '  -0.93%  '

$ws.Range('D4').Value = '''0.999'

$ws.Range('D5').Value = '''560.34'
$ws.Range('E5').Value = '  +0.04%  '

$ws.Range('D6').Value = '''161.90'
$ws.Range('E6').Value = '  +0.12%  '

$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('D8').Value = '''0.515'
$ws.Range('E8').Value = '  +1.82%  '

$ws.Range('E9').Value = '  +11.52%  '

$ws.Range('E10').Value = '  -1.64%  '

$ws.Range('D11').Value = '''0.331'
$ws.Range('E11').Value = '  +0.00%  '

$ws.Range('E12').Value = '  -5.12%  '

$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = '''0.0000176'
$ws.Range('E13').Value = '  +4.79%  '

$ws.Range('B14').Value = 'WrappedBTC'
$ws.Range('C14').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D14').Value = '68.743.64'
$ws.Range('E14').Value = '  +0.09%  '

$ws.Range('D15').Value = '2.872.37'
$ws.Range('E15').Value = '  -1.15%  '

$ws.Range('D16').Value = '''23.16'
$ws.Range('E16').Value = '  -1.52%  '

$ws.Range('D17').Value = '2.425.15'
$ws.Range('E17').Value = '  -1.69%  '

$ws.Range('D18').Value = '''10.49'
$ws.Range('E18').Value = '  -1.46%  '

$ws.Range('D19').Value = '''338.47'
$ws.Range('E19').Value = '  +1.23%  '

$ws.Range('D20').Value = '''6.93'
$ws.Range('E20').Value = '  +0.44%  '

$ws.Range('E21').Value = '  +1.22%  '

$ws.Range('D22').Value = '''1.93'
$ws.Range('E22').Value = '  +2.68%  '

$ws.Range('D24').Value = '''66.96'
$ws.Range('E24').Value = '  +0.60%  '

$ws.Range('D25').Value = '''3.69'
$ws.Range('E25').Value = '  +1.60%  '

$ws.Range('D26').Value = '2.553.57'
$ws.Range('E26').Value = '  -1.67%  '

$ws.Range('E27').Value = '  -0.07%  '

$ws.Range('D28').Value = '''8.21'
$ws.Range('E28').Value = '  +0.70%  '

$ws.Range('D29').Value = '0.0₃0816'
$ws.Range('E29').Value = '  +0.48%  '

$ws.Range('D30').Value = '''7.12'
$ws.Range('E30').Value = '  -0.50%  '

$ws.Range('D32').Value = '''427.55'
$ws.Range('E32').Value = '  -0.28%  '

$ws.Range('E33').Value = '  +1.98%  '

$ws.Range('E34').Value = '  +0.04%  '

$ws.Range('D35').Value = '''159.61'
$ws.Range('E35').Value = '  +0.47%  '

$ws.Range('D36').Value = '''19.03'
$ws.Range('E36').Value = '  +0.02%  '

$ws.Range('E37').Value = '  +0.06%  '

$ws.Range('D38').Value = '''17.95'
$ws.Range('E38').Value = '  +1.26%  '

$ws.Range('E39').Value = '  -3.01%  '

$ws.Range('D40').Value = '''0.298'
$ws.Range('E40').Value = '  -0.69%  '

$ws.Range('D41').Value = '''1.51'
$ws.Range('E41').Value = '  +3.28%  '

$ws.Range('D42').Value = '''4.33'
$ws.Range('E42').Value = '  -1.64%  '

$ws.Range('D43').Value = '''1.08'
$ws.Range('E43').Value = '  +0.85%  '

$ws.Range('D44').Value = '''2.04'
$ws.Range('E44').Value = '  -0.74%  '

$ws.Range('E45').Value = '  +0.09%  '

$ws.Range('D46').Value = '''130.64'
$ws.Range('E46').Value = '  +0.29%  '

$ws.Range('D47').Value = '''0.0716'
$ws.Range('E47').Value = '  +0.53%  '

$ws.Range('D48').Value = '''0.481'
$ws.Range('E48').Value = '  -0.25%  '

$ws.Range('D49').Value = '''0.557'
$ws.Range('E49').Value = '  -0.25%  '

$ws.Range('D50').Value = '''0.0924'
$ws.Range('E50').Value = '  +1.82%  '

$ws.Range('E51').Value = '  +1.06%  '
